$d = $word.ActiveDocument

# Locate the empty paragraph right after item "7." (the first of the two
# trailing empty paragraphs at the end of the document) and insert the new
# "8. ..." sentence (text + inline regression-splines equation + text) into
# it, preserving its existing paragraph mark / pPr.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -eq 0) {
        $target = $p
        break
    }
}

$insertRange = $d.Range($target.Range.Start, $target.Range.Start)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="David" w:hAnsi="David" w:cs="David" w:hint="cs"/>
                <w:i/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:rtl/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">8. נממש כעת פונקציות בסיס מבוססות על </w:t>
            </w:r>
            <m:oMath>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="David" w:hint="cs"/>
                  <w:color w:val="000000" w:themeColor="text1"/>
                  <w:sz w:val="24"/>
                  <w:szCs w:val="24"/>
                  <w:lang w:val="en-US"/>
                </w:rPr>
                <m:t xml:space="preserve"> </m:t>
              </m:r>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="David"/>
                  <w:color w:val="000000" w:themeColor="text1"/>
                  <w:sz w:val="24"/>
                  <w:szCs w:val="24"/>
                  <w:lang w:val="en-US"/>
                </w:rPr>
                <m:t>regression splines</m:t>
              </m:r>
            </m:oMath>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="David" w:hAnsi="David" w:cs="David" w:hint="cs"/>
                <w:i/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:rtl/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">. </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($xml)
